# Apply cryptos list update (generated from commit diff)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "38.128.93"
$ws.Range("E2").Value = "  +2.89%  "
$ws.Range("D3").Value = "2.061.79"
$ws.Range("E3").Value = "  +2.18%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "230.33"
$ws.Range("E5").Value = "  +1.64%  "
$ws.Range("E6").Value = "  +1.56%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.12"
$ws.Range("E7").Value = "  +6.46%  "
$ws.Range("E8").Value = "  +0.03%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.387"
$ws.Range("E9").Value = "  +1.87%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0809"
$ws.Range("E10").Value = "  +2.27%  "
$ws.Range("E11").Value = "  -0.64%  "
$ws.Range("D12").Value = "2.369.56"
$ws.Range("E12").Value = "  +2.27%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "14.62"
$ws.Range("E13").Value = "  +2.47%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "20.74"
$ws.Range("E14").Value = "  +2.03%  "
$ws.Range("E15").Value = "  +1.39%  "
$ws.Range("E16").Value = "  +3.04%  "
$ws.Range("D17").Value = "2.064.07"
$ws.Range("E17").Value = "  +1.95%  "
$ws.Range("D18").Value = "38.075.80"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.20"
$ws.Range("E19").Value = "  +2.17%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.86"
$ws.Range("E20").Value = "  +1.43%  "
$ws.Range("D21").Value = "0.0₃0833"
$ws.Range("E21").Value = "  +1.43%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "224.73"
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("E24").Value = "  +0.75%  "
$ws.Range("E25").Value = "  +2.56%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.28"
$ws.Range("E26").Value = "  +0.94%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "165.64"
$ws.Range("E27").Value = "  +0.36%  "
$ws.Range("E28").Value = "  +6.78%  "
$ws.Range("E29").Value = "  +1.50%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "19.07"
$ws.Range("E30").Value = "  +1.80%  "
$ws.Range("E31").Value = "  +1.63%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.56"
$ws.Range("E32").Value = "  +0.95%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.60"
$ws.Range("E33").Value = "  +3.73%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0614"
$ws.Range("E34").Value = "  +0.10%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.98"
$ws.Range("E35").Value = "  +7.14%  "
$ws.Range("E36").Value = "  +1.32%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "6.03"
$ws.Range("E37").Value = "  +11.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "3.30"
$ws.Range("E38").Value = "  +5.28%  "
$ws.Range("E39").Value = "  +0.04%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "98.45"
$ws.Range("E40").Value = "  +3.39%  "
$ws.Range("E41").Value = "  +0.77%  "
$ws.Range("D42").Value = "1.478.63"
$ws.Range("E42").Value = "  +0.09%  "
$ws.Range("B43").Value = "Cronos"
$ws.Range("C43").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.0945"
$ws.Range("E43").Value = "  +2.46%  "
$ws.Range("B44").Value = "HuobiToken"
$ws.Range("C44").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.87"
$ws.Range("E44").Value = "  +4.42%  "
$ws.Range("B45").Value = "InjectiveProtocol"
$ws.Range("C45").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "16.84"
$ws.Range("E45").Value = "  +1.66%  "
$ws.Range("E46").Value = "  -0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "4.09"
$ws.Range("E47").Value = "  +15.46%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.02"
$ws.Range("E48").Value = "  +1.22%  "
$ws.Range("E49").Value = "  +2.27%  "
$ws.Range("E50").Value = "  -1.70%  "
$ws.Range("D51").Value = "2.256.50"
$ws.Range("E51").Value = "  +2.34%  "
